{"js": "// The discussion's closing paragraph originally opened with a long\n// passage about the iron/brass rods and literature-value comparisons\n// before finally reaching the \"Man sollte daher die Messung erneut...\"\n// conclusion. The final edit drops that whole lead-in, so the\n// paragraph now starts directly with the concluding recommendation.\n\nconst body = context.document.body;\n\n// Anchor 1: the very first word of the text that must be removed.\nconst startResults = body.search(\"Die Zuordnung durch Farbe legt nahe\", { matchCase: true, matchWholeWord: false });\nstartResults.load(\"items\");\n\n// Anchor 2: the text that must remain - deletion stops right before it.\nconst keepResults = body.search(\"Man sollte daher die Messung erneut\", { matchCase: true, matchWholeWord: false });\nkeepResults.load(\"items\");\n\nawait context.sync();\n\nif (startResults.items.length === 0 || keepResults.items.length === 0) {\n  throw new Error(\"Could not locate the expected anchor text in the document.\");\n}\n\nconst deleteStart = startResults.items[0].getRange(\"Start\");\nconst deleteEnd = keepResults.items[0].getRange(\"Start\");\n\nconst deleteRange = deleteStart.expandTo(deleteEnd);\ndeleteRange.delete();\n\nawait context.sync();\n", "ps1": "# The discussion's closing paragraph originally opened with a long\n# passage about the iron/brass rods and literature-value comparisons\n# before finally reaching the \"Man sollte daher die Messung erneut...\"\n# conclusion. The final edit drops that whole lead-in, so the\n# paragraph now starts directly with the concluding recommendation.\n\n$d = $word.ActiveDocument\n\n# Anchor 1: the very first word of the text that must be removed.\n$startRange = $d.Content.Duplicate\n$fStart = $startRange.Find\n$fStart.ClearFormatting()\n$fStart.Forward = $true\n$fStart.MatchCase = $true\n$fStart.MatchWholeWord = $false\n$fStart.MatchWildcards = $false\n$foundStart = $fStart.Execute(\"Die Zuordnung durch Farbe legt nahe\")\n\n# Anchor 2: the text that must remain - deletion stops right before it.\n$endRange = $d.Content.Duplicate\n$fEnd = $endRange.Find\n$fEnd.ClearFormatting()\n$fEnd.Forward = $true\n$fEnd.MatchCase = $true\n$fEnd.MatchWholeWord = $false\n$fEnd.MatchWildcards = $false\n$foundEnd = $fEnd.Execute(\"Man sollte daher die Messung erneut\")\n\nif (-not $foundStart -or -not $foundEnd) {\n    throw \"Could not locate the expected anchor text in the document.\"\n}\n\n$deleteRange = $d.Range($startRange.Start, $endRange.Start)\n$deleteRange.Delete()\n"}
